$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2244.4614
$ws.Range("I129").Value = 1053.7778
$ws.Range("J129").Value = 4923.5
$ws.Range("K129").Value = 3161.3334
$ws.Range("L129").Value = 14770.5
$ws.Range("M129").Value = 1838.6666
$ws.Range("N129").Value = -24770.5
$ws.Range("H137").Value = 27779474
$ws.Range("I137").Value = 1129.0358
$ws.Range("K137").Value = 3387.1074
$ws.Range("M137").Value = -837.1074000000003

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1472037
$ws.Range("I2").Value = 1469.3334
$ws.Range("K2").Value = 1469.3334
$ws.Range("M2").Value = -1356.3334
$ws.Range("H21").Value = 17200
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9626
$ws.Range("H36").Value = 2575.25
$ws.Range("I36").Value = 2600.3333
$ws.Range("K36").Value = 2600.3333
$ws.Range("M36").Value = -2254.3333
$ws.Range("H61").Value = 1868.1333
$ws.Range("I61").Value = 1239.6154
$ws.Range("J61").Value = 5953.5
$ws.Range("K61").Value = 1239.6154
$ws.Range("L61").Value = 5953.5
$ws.Range("M61").Value = -1027.6154
$ws.Range("N61").Value = -6377.5
$ws.Range("H74").Value = 3155.6042
$ws.Range("I74").Value = 519.4231
$ws.Range("K74").Value = 519.4231
$ws.Range("M74").Value = 354.5769
$ws.Range("H77").Value = 3155.6042
$ws.Range("I77").Value = 519.4231
$ws.Range("K77").Value = 2597.1155
$ws.Range("M77").Value = 1770.8845
$ws.Range("H110").Value = 851.0540999999999
$ws.Range("I110").Value = 735.92
$ws.Range("J110").Value = 1090.9166
$ws.Range("K110").Value = 735.92
$ws.Range("L110").Value = 1090.9166
$ws.Range("M110").Value = 1309.08
$ws.Range("N110").Value = -5180.9166
$ws.Range("H116").Value = 1472037
$ws.Range("I116").Value = 1469.3334
$ws.Range("K116").Value = 1469.3334
$ws.Range("M116").Value = 824.6666
$ws.Range("H132").Value = 1962.4412
$ws.Range("I132").Value = 1775
$ws.Range("J132").Value = 2685.4285
$ws.Range("K132").Value = 5325
$ws.Range("L132").Value = 8056.2855
$ws.Range("M132").Value = -2795
$ws.Range("N132").Value = -13116.2855
$ws.Range("H136").Value = 1868.1333
$ws.Range("I136").Value = 1239.6154
$ws.Range("J136").Value = 5953.5
$ws.Range("K136").Value = 3718.8462
$ws.Range("L136").Value = 17860.5
$ws.Range("M136").Value = -1168.8462
$ws.Range("N136").Value = -22960.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1472037
$ws.Range("I3").Value = 1469.3334
$ws.Range("K3").Value = 1469.3334
$ws.Range("M3").Value = -1355.3334
$ws.Range("H94").Value = 364.9524
$ws.Range("I94").Value = 292.8421
$ws.Range("K94").Value = 292.8421
$ws.Range("M94").Value = 158.1579
$ws.Range("H134").Value = 40159.367
$ws.Range("I134").Value = 53091.684
$ws.Range("J134").Value = 4595.5
$ws.Range("K134").Value = 159275.052
$ws.Range("L134").Value = 13786.5
$ws.Range("M134").Value = -156740.052
$ws.Range("N134").Value = -18856.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 6555.5
$ws.Range("I17").Value = 6555.5
$ws.Range("K17").Value = 6555.5
$ws.Range("M17").Value = -6381.5
$ws.Range("H31").Value = 1381.3636
$ws.Range("I31").Value = 1256.6666
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1256.6666
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -961.6666
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 1381.3636
$ws.Range("I34").Value = 1256.6666
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1256.6666
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1054.6666
$ws.Range("N34").Value = -4404
$ws.Range("H58").Value = 2187.5
$ws.Range("I58").Value = 1741.5
$ws.Range("J58").Value = 3748.5
$ws.Range("K58").Value = 1741.5
$ws.Range("L58").Value = 3748.5
$ws.Range("M58").Value = -1538.5
$ws.Range("N58").Value = -4154.5
$ws.Range("H132").Value = 2264.2144
$ws.Range("I132").Value = 1614.5714
$ws.Range("J132").Value = 4213.143
$ws.Range("K132").Value = 4843.7142
$ws.Range("L132").Value = 12639.429
$ws.Range("M132").Value = -2313.7142
$ws.Range("N132").Value = -17699.429
$ws.Range("H134").Value = 2397.7727
$ws.Range("I134").Value = 2313.2104
$ws.Range("J134").Value = 2933.3333
$ws.Range("K134").Value = 6939.6312
$ws.Range("L134").Value = 8799.999899999999
$ws.Range("M134").Value = -4404.6312
$ws.Range("N134").Value = -13869.9999
$ws.Range("H136").Value = 2187.5
$ws.Range("I136").Value = 1741.5
$ws.Range("J136").Value = 3748.5
$ws.Range("K136").Value = 5224.5
$ws.Range("L136").Value = 11245.5
$ws.Range("M136").Value = -2674.5
$ws.Range("N136").Value = -16345.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 511.25
$ws.Range("I15").Value = 152.5
$ws.Range("J15").Value = 870
$ws.Range("K15").Value = 457.5
$ws.Range("L15").Value = 2610
$ws.Range("M15").Value = -317.5
$ws.Range("N15").Value = -2890
$ws.Range("H44").Value = 1882.1818
$ws.Range("I44").Value = 500
$ws.Range("J44").Value = 2400.5
$ws.Range("K44").Value = 1500
$ws.Range("L44").Value = 7201.5
$ws.Range("M44").Value = -1102
$ws.Range("N44").Value = -7997.5
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10372
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 27000
$ws.Range("N65").Value = -33864
$ws.Range("H131").Value = 2275751.2
$ws.Range("J131").Value = 3228333.2
$ws.Range("L131").Value = 9684999.600000001
$ws.Range("N131").Value = -9695079.600000001
$ws.Range("H134").Value = 4209.9414
$ws.Range("I134").Value = 2256.9
$ws.Range("K134").Value = 6770.700000000001
$ws.Range("M134").Value = -1700.700000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3840.9473
$ws.Range("I132").Value = 3755
$ws.Range("J132").Value = 3988.2856
$ws.Range("K132").Value = 11265
$ws.Range("L132").Value = 11964.8568
$ws.Range("M132").Value = -8735
$ws.Range("N132").Value = -17024.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2005
$ws.Range("I31").Value = 1015
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1015
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -767
$ws.Range("N31").Value = -2996
$ws.Range("H97").Value = 20200
$ws.Range("J97").Value = 20200
$ws.Range("L97").Value = 20200
$ws.Range("N97").Value = -22182
$ws.Range("H132").Value = 2042.6857
$ws.Range("I132").Value = 1351.6316
$ws.Range("J132").Value = 2863.3125
$ws.Range("K132").Value = 4054.8948
$ws.Range("L132").Value = 8589.9375
$ws.Range("M132").Value = -1524.8948
$ws.Range("N132").Value = -13649.9375
$ws.Range("H136").Value = 1794.9131
$ws.Range("I136").Value = 1326.5555
$ws.Range("J136").Value = 3481
$ws.Range("K136").Value = 3979.6665
$ws.Range("L136").Value = 10443
$ws.Range("M136").Value = -1429.6665
$ws.Range("N136").Value = -15543

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 11635
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 14952.5
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 14952.5
$ws.Range("M32").Value = -4683
$ws.Range("N32").Value = -15586.5
$ws.Range("H113").Value = 461.83334
$ws.Range("I113").Value = 501.54544
$ws.Range("K113").Value = 1504.63632
$ws.Range("M113").Value = 665.3636799999999
$ws.Range("H132").Value = 2115.4666
$ws.Range("I132").Value = 1123.4
$ws.Range("J132").Value = 4099.6
$ws.Range("K132").Value = 3370.2
$ws.Range("L132").Value = 12298.8
$ws.Range("M132").Value = -840.2000000000003
$ws.Range("N132").Value = -17358.8
$ws.Range("H136").Value = 15413.929
$ws.Range("I136").Value = 15413.929
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 46241.787
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -43691.787
$ws.Range("N136").ClearContents()
